$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Clear()

$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Rows.Item(1).AutoFit()

$ws.Range("B2").Value = 'LOQ4203'
$ws.Range("C2").Value = 'LOQ4203'
$ws.Rows.Item(2).AutoFit()

$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Sistemas Produtivos I'
$ws.Range("C3").Value = ' Sistemas Produtivos I'
$ws.Rows.Item(3).AutoFit()

$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Productive Systems I'
$ws.Range("C4").Value = 'Productive Systems I'
$ws.Rows.Item(4).AutoFit()

$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '2'
$ws.Range("C5").Value = '2'
$ws.Rows.Item(5).AutoFit()

$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '0'
$ws.Range("C6").Value = '0'
$ws.Rows.Item(6).AutoFit()

$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '30 h'
$ws.Range("C7").Value = '30 h'
$ws.Rows.Item(7).AutoFit()

$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2018'
$ws.Range("C8").Value = '01/01/2018'
$ws.Rows.Item(8).AutoFit()

$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EP-2'
$ws.Range("C9").Value = 'EP-2'
$ws.Rows.Item(9).AutoFit()

$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = 'Introduzir os alunos nos conceitos técnicos fundamentais de um curso de Engenharia de Produção, tendo em vista a sua formação generalista voltada para os mais diversos tipos de sistemas de produção.'
$ws.Range("C10").Value = 'Introduzir os alunos nos conceitos técnicos fundamentais de um curso de Engenharia de Produção, tendo em vista a sua formação generalista voltada para os mais diversos tipos de sistemas de produção.'
$ws.Rows.Item(10).RowHeight = 60

$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = 'Introduce students to the fundamental technical concepts of a Industrial Engineering course, with a view to their general training aimed at the most diverse types of production systems.'
$ws.Range("C11").Value = 'Introduce students to the fundamental technical concepts of a Industrial Engineering course, with a view to their general training aimed at the most diverse types of production systems.'
$ws.Rows.Item(11).RowHeight = 60

$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Rows.Item(12).AutoFit()

$ws.Range("B13").Value = '5840535 - Messias Borges Silva'
$ws.Range("C13").Value = '5840535 - Messias Borges Silva'
$ws.Rows.Item(13).AutoFit()

$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = '1 – Introdução aos Sistemas Produtivos2 – Papel Estratégico da Produção3 – Estratégia de Produção4 – Projeto em Gestão de Produção5 – Projeto de Produtos e Serviços6 – Projeto da Rede de Operações Produtivas7 – Arranjo Físico e Fluxo'
$ws.Range("C14").Value = '1 – Introdução aos Sistemas Produtivos2 – Papel Estratégico da Produção3 – Estratégia de Produção4 – Projeto em Gestão de Produção5 – Projeto de Produtos e Serviços6 – Projeto da Rede de Operações Produtivas7 – Arranjo Físico e Fluxo'
$ws.Rows.Item(14).RowHeight = 60

$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = '1 - Introduction to Productive Systems2 - Strategic Role of Production3 - Production Strategy4 - Project in Production Management5 - Product and Service Project6 - Production Operations Network Project7 - Physical Arrangement and Flow'
$ws.Range("C15").Value = '1 - Introduction to Productive Systems2 - Strategic Role of Production3 - Production Strategy4 - Project in Production Management5 - Product and Service Project6 - Production Operations Network Project7 - Physical Arrangement and Flow'
$ws.Rows.Item(15).RowHeight = 60

$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = '1 – Introdução aos Sistemas ProdutivosProdução na Organização. Inputs, Processos de Transformação e Outputs. Tipos de Operações de Produção. Atividades da administração da produção.2 – Papel Estratégico da ProduçãoPapel da função produção. Objetivos de desempenho. 3 – Estratégia de ProduçãoIntrodução. Prioridade de objetivos de desempenho. Áreas de decisão da estratégia de operações.4 – Projeto em Gestão de ProduçãoDefinição de projeto. Principais aspectos de um projeto. Tipos de processos em manufatura e serviços. 5 – Projeto de Produtos e ServiçosIntrodução. Geração de conceito. Triagem de conceito. Avaliação e melhoria do projeto. Protótipo e projeto final.6 – Projeto da Rede de Operações ProdutivasPerspectiva da rede. Integração Vertical. Localização da capacidade. Gestão da capacidade produtiva de longo prazo.7 – Arranjo Físico e FluxoProcedimento de Arranjo Físico. Tipos básicos de arranjo físico. Projeto de arranjo físico.'
$ws.Range("C16").Value = '1 – Introdução aos Sistemas ProdutivosProdução na Organização. Inputs, Processos de Transformação e Outputs. Tipos de Operações de Produção. Atividades da administração da produção.2 – Papel Estratégico da ProduçãoPapel da função produção. Objetivos de desempenho. 3 – Estratégia de ProduçãoIntrodução. Prioridade de objetivos de desempenho. Áreas de decisão da estratégia de operações.4 – Projeto em Gestão de ProduçãoDefinição de projeto. Principais aspectos de um projeto. Tipos de processos em manufatura e serviços. 5 – Projeto de Produtos e ServiçosIntrodução. Geração de conceito. Triagem de conceito. Avaliação e melhoria do projeto. Protótipo e projeto final.6 – Projeto da Rede de Operações ProdutivasPerspectiva da rede. Integração Vertical. Localização da capacidade. Gestão da capacidade produtiva de longo prazo.7 – Arranjo Físico e FluxoProcedimento de Arranjo Físico. Tipos básicos de arranjo físico. Projeto de arranjo físico.'
$ws.Rows.Item(16).RowHeight = 120

$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = '1 - Introduction to Productive SystemsProduction in the Organization. Inputs, Transformation Processes and Outputs. Types of Production Operations. Production management activities.2 - Strategic Role of ProductionRole of production function. Performance Objectives.3 - Production StrategyIntroduction. Priority of performance goals. Operations strategy decision areas.4 - Project in Production ManagementDefinition of project. Main aspects of a project. Types of processes in manufacturing and services.5 - Product and Service ProjectIntroduction. Concept generation. Concept screening. Evaluation and improvement of the project. Prototype and final design.6 - Production Operations Network ProjectNetwork perspective. Vertical integration. Location of capacity. Management of long-term productive capacity.7 – Layout and FlowLayout and Physical Arrangement Procedure. Basic types of physical arrangement. Design of layout and physical arrangement.'
$ws.Range("C17").Value = '1 - Introduction to Productive SystemsProduction in the Organization. Inputs, Transformation Processes and Outputs. Types of Production Operations. Production management activities.2 - Strategic Role of ProductionRole of production function. Performance Objectives.3 - Production StrategyIntroduction. Priority of performance goals. Operations strategy decision areas.4 - Project in Production ManagementDefinition of project. Main aspects of a project. Types of processes in manufacturing and services.5 - Product and Service ProjectIntroduction. Concept generation. Concept screening. Evaluation and improvement of the project. Prototype and final design.6 - Production Operations Network ProjectNetwork perspective. Vertical integration. Location of capacity. Management of long-term productive capacity.7 – Layout and FlowLayout and Physical Arrangement Procedure. Basic types of physical arrangement. Design of layout and physical arrangement.'
$ws.Rows.Item(17).RowHeight = 120

$ws.Range("A18").Value = 'Avaliação:'
$ws.Rows.Item(18).AutoFit()

$ws.Range("A19").Value = 'Método:'
$ws.Range("B19").Value = 'Aulas Expositivas; trabalhos e seminários.'
$ws.Range("C19").Value = 'Aulas Expositivas; trabalhos e seminários.'
$ws.Rows.Item(19).RowHeight = 60

$ws.Range("A20").Value = 'Critério:'
$ws.Range("B20").Value = 'MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.'
$ws.Range("C20").Value = 'MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.'
$ws.Rows.Item(20).RowHeight = 60

$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.'
$ws.Rows.Item(21).RowHeight = 60

$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'SLACK, N. et al. Administração da produção. São Paulo: Atlas, 2002. 
Textos complementares serão usados durante o curso.'
$ws.Range("C22").Value = 'SLACK, N. et al. Administração da produção. São Paulo: Atlas, 2002. 
Textos complementares serão usados durante o curso.'
$ws.Rows.Item(22).RowHeight = 120

$ws.Range("A23").Value = 'Requisitos:'
$ws.Rows.Item(23).AutoFit()

$ws.Range("B24").Value = 'LOQ4201 -  Introdução à Engenharia de Produção  (Requisito fraco)
'
$ws.Range("C24").Value = 'LOQ4201 -  Introdução à Engenharia de Produção  (Requisito fraco)
'
$ws.Rows.Item(24).RowHeight = 30

